# Atualização de bases das ligas, do dia: 24-02-2024 às 21:58
#
# The underlying data rows got re-paired / re-ordered: for a number of row
# pairs, the entire record (columns B..AC -- i.e. everything except the
# leading sequential index in column A) was swapped between two adjacent
# rows. Column A (the running id/index) stays put for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple is a pair of row numbers whose B:AC content must be swapped.
$rowPairs = @(
    @(68, 69),
    @(79, 80),
    @(98, 99),
    @(124, 125),
    @(141, 142),
    @(236, 237),
    @(298, 299),
    @(380, 381),
    @(386, 387),
    @(395, 396)
)

# Column B is index 2, column AC is index 29.
$colStart = 2
$colEnd = 29

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    for ($c = $colStart; $c -le $colEnd; $c++) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}
